$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the "Supplementary Table" style row (row 7) onto
# row 25 (fixing its previously mismatched style) and onto the brand new
# row 26.
$ws.Range("A7:E7").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Range("A26:E26").PasteSpecial(-4122)

# Row 25: Supplementary Table S3.4 (existing row, now gains D/E columns)
$ws.Range("A25").Value = "Supplementary Table S3.4"
$ws.Range("B25").Value = "Online Supplementary Material"
$ws.Range("C25").Value = "List of keywords and profiles for phototransduction gene families."
$ws.Range("D25").Value = "Prepared. To be added in github."
$ws.Range("E25").Value = ""

# Row 26: new Supplementary Table S3.5 entry
$ws.Range("A26").Value = "Supplementary Table S3.5"
$ws.Range("B26").Value = "Online Supplementary Material"
$ws.Range("C26").Value = "List of models for gene trees"
$ws.Range("D26").Value = "Prepared. To be added in github."
$ws.Range("E26").Value = ""

$null = $ws.Range("A27").Select()
